# Update "想去人数" (want-to-go count) values in column F
# for the "展览" (Exhibition) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Updates that are identical between the two sheets (row -> new value)
$commonUpdates = @{
    2  = 1043
    8  = 1650
    9  = 6062
    11 = 342
    12 = 274
    17 = 257
    23 = 250
    25 = 158
    26 = 4
    28 = 376
    29 = 68
    32 = 39
    33 = 53
    35 = 57
}

foreach ($ws in @($wsExhibition, $wsAllTypes)) {
    foreach ($row in $commonUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $commonUpdates[$row]
    }
}

# Row 16 diverges between the two sheets after the edit
$wsExhibition.Cells.Item(16, 6).Value = 5359
$wsAllTypes.Cells.Item(16, 6).Value = 5361
